$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Copy formatting (styles) of the previous changelog entry row (20) into the
# new row (21) so that the new row matches the look of the existing table.
$ws.Range("A20:C20").Copy()
$ws.Range("A21:C21").PasteSpecial(-4122)

# Fill in the new changelog entry: date, version, changes
$ws.Range("A21").Value = 44369
$ws.Range("B21").Value = "1.5.2"
$ws.Range("C21").Value = "Bugs:`n- Wrong representation of links (color is standard blue instead of pink) fixed"

# Match the row height used by other multi-line changelog rows
$ws.Rows.Item(21).RowHeight = 30

# Move active selection to the cell below the newly added row, as in the source file
$ws.Range("C22").Select()
